$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the last existing data row (A61) onto the
# new rows so the added dates reuse the existing date style instead of
# creating a new one.
$ws.Range("A61").Copy()
$ws.Range("A62:A64").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Append three new rows of data (Date, Error Count) using clean date serials
$ws.Range("A62").Value = 46043
$ws.Range("B62").Value = 6

$ws.Range("A63").Value = 46044
$ws.Range("B63").Value = 5

$ws.Range("A64").Value = 46042
$ws.Range("B64").Value = 2

# Update the selected cell/range to reflect the edit location
$ws.Range("E62").Select()
